$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff
$ws.Range("C2").Value = 12.2
$ws.Range("B3").Value = 4.5999999999999996
$ws.Range("B4").Value = 0.65

# Update column widths
$ws.Columns.Item(1).ColumnWidth = 27
$ws.Columns.Item(2).ColumnWidth = 8.43
$ws.Columns.Item(3).ColumnWidth = 27.25

# Update selection
$ws.Range("B4").Select()
